# Updated cryptos list on Sat Sep  9 14:30:11 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.069.74"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.645.80"
$ws.Range("E3").Value = "  +0.85%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.94%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.93"
$ws.Range("E5").Value = "  +1.04%  "

# Row 6: XRP
$ws.Range("E6").Value = "  +1.06%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.96%  "

# Row 8: Cardano
$ws.Range("E8").Value = "  +0.45%  "

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0640"
$ws.Range("E9").Value = "  +1.30%  "

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("E10").Value = "  +0.12%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("E11").Value = "  +0.87%  "

# Row 12: Polkadot
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.874.34"
$ws.Range("E12").Value = "  +0.90%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.30"
$ws.Range("E13").Value = "  +1.59%  "

# Row 14: WrappedEther
$ws.Range("D14").Value = "1.669.24"
$ws.Range("E14").Value = "  +2.46%  "

# Row 15: Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.545"
$ws.Range("E15").Value = "  +0.03%  "

# Row 16: ShibaInu
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  +1.23%  "

# Row 17: Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.44"
$ws.Range("E17").Value = "  +1.01%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "26.158.79"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19: Dai
$ws.Range("E19").Value = "  +0.94%  "

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.46"
$ws.Range("E20").Value = "  +0.00%  "

# Row 21: Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.34"
$ws.Range("E21").Value = "  -0.91%  "

# Row 22: Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23: Chainlink
$ws.Range("E23").Value = "  -0.50%  "

# Row 24: Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.81"
$ws.Range("E24").Value = "  +1.25%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.47"
$ws.Range("E25").Value = "  +1.56%  "

# Row 26: BinanceUSD
$ws.Range("E26").Value = "  +1.24%  "

# Row 27: Stellar
$ws.Range("E27").Value = "  +3.91%  "

# Row 28: Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.91"
$ws.Range("E28").Value = "  +0.70%  "

# Row 29: EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +0.57%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  +1.27%  "

# Row 31: Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0498"
$ws.Range("E31").Value = "  -0.32%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("E32").Value = "  +1.55%  "

# Row 33: Filecoin
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.29"
$ws.Range("E33").Value = "  -1.01%  "

# Row 34: LidoDAOToken
$ws.Range("E34").Value = "  -3.13%  "

# Row 35: HuobiToken
$ws.Range("E35").Value = "  +2.01%  "

# Row 36: ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.906"
$ws.Range("E36").Value = "  +0.48%  "

# Row 37: Maker
$ws.Range("D37").Value = "1.135.00"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38: ImmutableX
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.541"
$ws.Range("E38").Value = "  -1.76%  "

# Row 39: MXToken
$ws.Range("E39").Value = "  +0.02%  "

# Row 40: VeChain
$ws.Range("E40").Value = "  +0.45%  "

# Row 41: FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.50"
$ws.Range("E41").Value = "  +1.03%  "

# Row 42: Quant
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.65"
$ws.Range("E42").Value = "  +0.48%  "

# Row 43: TrustWalletToken
$ws.Range("E43").Value = "  -0.81%  "

# Row 44: RocketPoolETH
$ws.Range("D44").Value = "1.784.15"

# Row 45: BabyDogeCoin
$ws.Range("D45").Value = "0.0₆0117"
$ws.Range("E45").Value = "  +5.52%  "

# Row 46: Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.72"
$ws.Range("E46").Value = "  +1.00%  "

# Row 47: Cronos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0531"
$ws.Range("E47").Value = "  +0.85%  "

# Row 48: RenderToken
$ws.Range("E48").Value = "  -0.29%  "

# Row 49: EnergySwap
$ws.Range("E49").Value = "  +1.65%  "

# Row 50: Mantle
$ws.Range("E50").Value = "  +0.54%  "

# Row 51: Algorand
$ws.Range("E51").Value = "  -0.28%  "
